$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.036306149394157
$ws.Range("D2").Value = 1.043600749504522
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.050888568610685
$ws.Range("I2").Value = 1.036444174204034
$ws.Range("J2").Value = 1.041415186115142
$ws.Range("K2").Value = 1.046374157013483
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.053641595100632
$ws.Range("N2").Value = 1.017683858959257

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.037605606681334
$ws.Range("D3").Value = 1.044640422231045
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.052132180743302
$ws.Range("I3").Value = 1.036758847152743
$ws.Range("J3").Value = 1.042356790447942
$ws.Range("K3").Value = 1.047224300296608
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.05469663410175
$ws.Range("N3").Value = 1.018000639594662

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.038445353041566
$ws.Range("D4").Value = 1.045311889249536
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.052936023679849
$ws.Range("I4").Value = 1.036960172815736
$ws.Range("J4").Value = 1.042964527428847
$ws.Range("K4").Value = 1.047772530998591
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.055377900122462
$ws.Range("N4").Value = 1.018204966033602

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.038798125645813
$ws.Range("D5").Value = 1.045593872120193
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.053273757049421
$ws.Range("I5").Value = 1.037044263104284
$ws.Range("J5").Value = 1.043219653070409
$ws.Range("K5").Value = 1.048002561989118
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.055663968352341
$ws.Range("N5").Value = 1.018290709702559

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.038857342778392
$ws.Range("D6").Value = 1.045641200651578
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.053330452198752
$ws.Range("I6").Value = 1.037058350182604
$ws.Range("J6").Value = 1.043262468344595
$ws.Range("K6").Value = 1.048041159161852
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.055711980841159
$ws.Range("N6").Value = 1.018305097368315

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.038450067810292
$ws.Range("D7").Value = 1.045315658301687
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.052940537280442
$ws.Range("I7").Value = 1.036961298582596
$ws.Range("J7").Value = 1.042967937868903
$ws.Range("K7").Value = 1.047775606430563
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.055381723897118
$ws.Range("N7").Value = 1.018206112353313

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.036745535072934
$ws.Range("D8").Value = 1.043952376315288
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.051309032489265
$ws.Range("I8").Value = 1.036550993863393
$ws.Range("J8").Value = 1.041733727000836
$ws.Range("K8").Value = 1.046661855211145
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.053998445033738
$ws.Range("N8").Value = 1.017791051952093

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.033733405248562
$ws.Range("D9").Value = 1.041540259833056
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.048427409870793
$ws.Range("I9").Value = 1.035810413236598
$ws.Range("J9").Value = 1.03954694497272
$ws.Range("K9").Value = 1.044684869862334
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.051549973739782
$ws.Range("N9").Value = 1.017054630644306

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.031719329537067
$ws.Range("D10").Value = 1.039925413492661
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.046501626415455
$ws.Range("I10").Value = 1.035304817899793
$ws.Range("J10").Value = 1.038080888954356
$ws.Range("K10").Value = 1.043357045405898
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.04991011879534
$ws.Range("N10").Value = 1.016560243449318

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.030845733871083
$ws.Range("D11").Value = 1.039224526564188
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.045666581234682
$ws.Range("I11").Value = 1.035083058697552
$ws.Range("J11").Value = 1.037444086464644
$ws.Range("K11").Value = 1.042779717281112
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.049198215237899
$ws.Range("N11").Value = 1.016345339616539

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.030521012936867
$ws.Range("D12").Value = 1.038963935229148
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.045356228803891
$ws.Range("I12").Value = 1.035000260494343
$ws.Range("J12").Value = 1.037207247153502
$ws.Range("K12").Value = 1.042564912678852
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.048933503246436
$ws.Range("N12").Value = 1.016265388779615

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.030590677106985
$ws.Range("D13").Value = 1.039019844351906
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.045422808639108
$ws.Range("I13").Value = 1.035018040332331
$ws.Range("J13").Value = 1.037258063734827
$ws.Range("K13").Value = 1.042611005299007
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.048990297563911
$ws.Range("N13").Value = 1.016282544225161

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.030818897030584
$ws.Range("D14").Value = 1.039202991136198
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.045640931090084
$ws.Range("I14").Value = 1.035076223289847
$ws.Range("J14").Value = 1.037424515444558
$ws.Range("K14").Value = 1.042761968811112
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.049176339780613
$ws.Range("N14").Value = 1.016338733432149

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.030959480411489
$ws.Range("D15").Value = 1.039315800596359
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.04577529963059
$ws.Range("I15").Value = 1.035112015113824
$ws.Range("J15").Value = 1.037527031646805
$ws.Range("K15").Value = 1.042854934722686
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.049290929412301
$ws.Range("N15").Value = 1.016373336728308

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.031777275426847
$ws.Range("D16").Value = 1.039971894063542
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.046557020663681
$ws.Range("I16").Value = 1.035319475525351
$ws.Range("J16").Value = 1.038123109205423
$ws.Range("K16").Value = 1.043395310558157
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.049957326515371
$ws.Range("N16").Value = 1.016574488304615

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.032289854764711
$ws.Range("D17").Value = 1.040383000777009
$ws.Range("E17").Value = 0.9894763578477731
$ws.Range("F17").Value = 1.047047058062477
$ws.Range("I17").Value = 1.035448850534623
$ws.Range("J17").Value = 1.038496477551197
$ws.Range("K17").Value = 1.043733637079272
$ws.Range("L17").Value = 0.9930127773692701
$ws.Range("M17").Value = 1.05037484589325
$ws.Range("N17").Value = 1.016700442105839

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.032588690437554
$ws.Range("D18").Value = 1.040622633594328
$ws.Range("E18").Value = 0.9897087662937551
$ws.Range("F18").Value = 1.047332776138751
$ws.Range("I18").Value = 1.035524039527922
$ws.Range("J18").Value = 1.038714065166084
$ws.Range("K18").Value = 1.043930748614367
$ws.Range("L18").Value = 0.9932001317071766
$ws.Range("M18").Value = 1.050618201021113
$ws.Range("N18").Value = 1.016773828734308

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.032690561546824
$ws.Range("D19").Value = 1.040704315356284
$ws.Range("E19").Value = 0.9897880325774039
$ws.Range("F19").Value = 1.047430179664756
$ws.Range("I19").Value = 1.035549630721334
$ws.Range("J19").Value = 1.038788224486721
$ws.Range("K19").Value = 1.043997919868064
$ws.Range("L19").Value = 0.993264023964098
$ws.Range("M19").Value = 1.050701148938352
$ws.Range("N19").Value = 1.016798838127506

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.032234874701314
$ws.Range("D20").Value = 1.040338909364006
$ws.Range("E20").Value = 0.9894336180360677
$ws.Range("F20").Value = 1.046994493348208
$ws.Range("I20").Value = 1.035434998095274
$ws.Range("J20").Value = 1.038456438528628
$ws.Range("K20").Value = 1.043697361513145
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.05033006833582
$ws.Range("N20").Value = 1.016686936748064

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.030751698303236
$ws.Range("D21").Value = 1.039149065919742
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.045576704451032
$ws.Range("I21").Value = 1.035059101655971
$ws.Range("J21").Value = 1.037375507943463
$ws.Range("K21").Value = 1.042717523759918
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.049121562723788
$ws.Range("N21").Value = 1.016322190598243

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.029817841226599
$ws.Range("D22").Value = 1.038399511929231
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.044684243955857
$ws.Range("I22").Value = 1.034820289709313
$ws.Range("J22").Value = 1.036694132004807
$ws.Range("K22").Value = 1.042099380986936
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.04836010948484
$ws.Range("N22").Value = 1.016092130732219

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.030313023862335
$ws.Range("D23").Value = 1.038797003394972
$ws.Range("E23").Value = 0.9879432794636464
$ws.Range("F23").Value = 1.04515745409187
$ws.Range("I23").Value = 1.034947123102665
$ws.Range("J23").Value = 1.037055509509446
$ws.Range("K23").Value = 1.042427268311095
$ws.Range("L23").Value = 0.9917760702887611
$ws.Range("M23").Value = 1.04876392479292
$ws.Range("N23").Value = 1.0162141593526

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.032259718267492
$ws.Range("D24").Value = 1.040358832866116
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.04701824542724
$ws.Range("I24").Value = 1.035441258260625
$ws.Range("J24").Value = 1.038474531032
$ws.Range("K24").Value = 1.043713753586221
$ws.Range("L24").Value = 0.9929938892766441
$ws.Range("M24").Value = 1.05035030193291
$ws.Range("N24").Value = 1.016693039485334

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.03451314967081
$ws.Range("D25").Value = 1.042165031448804
$ws.Range("E25").Value = 0.9912096547607051
$ws.Range("F25").Value = 1.049173192676858
$ws.Range("I25").Value = 1.036003959583983
$ws.Range("J25").Value = 1.040113713500064
$ws.Range("K25").Value = 1.04519769069294
$ws.Range("L25").Value = 0.9944092447426416
$ws.Range("M25").Value = 1.052184278521951
$ws.Range("N25").Value = 1.017245615294967

